# The train-status checker now keeps running (instead of a single one-shot
# check) and refreshes each train's "Delay" cell with the latest status it
# read back from the incoming status e-mails, including the new "departed
# on time from București Nord" replies. Push the refreshed statuses onto
# the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value  = "+6 min întârziere"
$ws.Range("B5").Value  = "+3 min întârziere"
$ws.Range("B6").Value  = "+6 min întârziere"
$ws.Range("B7").Value  = "la timp"
$ws.Range("B8").Value  = "la timp"
$ws.Range("B9").Value  = "la timp"
$ws.Range("B10").Value = "la timp"
$ws.Range("B11").Value = "sosește cu 35 min întârziere la Constanța*"
$ws.Range("B12").Value = "sosește la timp la Constanța*"
$ws.Range("B13").Value = "pleacă la timp din București Nord*"
$ws.Range("B14").Value = "pleacă la timp din București Nord*"
$ws.Range("B15").Value = "pleacă la timp din București Nord*"
$ws.Range("B16").Value = "pleacă la timp din București Nord*"
